$wb = $excel.ActiveWorkbook

# "展览" sheet — filtered "exhibition" listing, rows 3-6 hold the four updated events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 211
$wsExhibit.Range("F4").Value = 2467
$wsExhibit.Range("F5").Value = 34
$wsExhibit.Range("F6").Value = 537

# "全部类型" sheet — combined listing, same four events sit two rows lower (rows 5-8)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 211
$wsAll.Range("F6").Value = 2467
$wsAll.Range("F7").Value = 34
$wsAll.Range("F8").Value = 537
